$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "2degrees, telecommunications" -> append " - " + hyperlink + "; " +
#    "National Sales Manager"
# ---------------------------------------------------------------------------
$idx2degrees = Find-ParagraphIndex $d "2degrees, telecommunications"
$p = $d.Paragraphs.Item($idx2degrees)
$r = $p.Range
$r.InsertAfter(" - QQQQQQQQQQQ")
$phEnd = $r.End - 1
$phStart = $phEnd - 11
$hlRange = $d.Range($phStart, $phEnd)
$null = $d.Hyperlinks.Add($hlRange, "mailto:cory.moran@2degreesmobile.co.nz", [System.Type]::Missing, [System.Type]::Missing, "cory.moran@2degreesmobile.co.nz")

$p = $d.Paragraphs.Item($idx2degrees)
$p.Range.InsertAfter("; ")
$p = $d.Paragraphs.Item($idx2degrees)
$p.Range.InsertAfter("National Sales Manager")

# ---------------------------------------------------------------------------
# 2) "Spark New Zealand " -> append "- " + hyperlink + "; group HR Director"
# ---------------------------------------------------------------------------
$idxSpark = Find-ParagraphIndex $d "Spark New Zealand"
$p = $d.Paragraphs.Item($idxSpark)
$r = $p.Range
$r.InsertAfter("- QQQQQQQQQQQQ")
$phEnd = $r.End - 1
$phStart = $phEnd - 12
$hlRange = $d.Range($phStart, $phEnd)
$null = $d.Hyperlinks.Add($hlRange, "mailto:joe.mccollum@spark.co.nz", [System.Type]::Missing, [System.Type]::Missing, "joe.mccollum@spark.co.nz")

$p = $d.Paragraphs.Item($idxSpark)
$p.Range.InsertAfter("; group HR Director")

# ---------------------------------------------------------------------------
# 3) Clean up "SLI Systems, eCommerce solutions, sea|rch engines" paragraph:
#    remove the proofErr spell-check run splits and the _GoBack bookmark that
#    currently sits in the middle of the text, collapsing everything into a
#    single plain run.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $gb = $d.Bookmarks.Item("_GoBack")
    $gb.Delete()
}

$idxSli = Find-ParagraphIndex $d "SLI Systems"
$p = $d.Paragraphs.Item($idxSli)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "TEMP_PLACEHOLDER_TEXT_ZZZ"
$p = $d.Paragraphs.Item($idxSli)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = "SLI Systems, eCommerce solutions, search engines"

# ---------------------------------------------------------------------------
# 4) Re-add the _GoBack bookmark right after "Warehouse group (and financial
#    services)" (its new home per the authors edit).
# ---------------------------------------------------------------------------
$idxWarehouse = Find-ParagraphIndex $d "Warehouse group"
$p = $d.Paragraphs.Item($idxWarehouse)
$r = $p.Range
$origEnd = $r.End
$r.InsertAfter("TEMP_BOOKMARK_ANCHOR")
$bmPos = $origEnd - 1
$bmRange = $d.Range($bmPos, $bmPos)
$null = $d.Bookmarks.Add("_GoBack", $bmRange)
$cleanupRange = $d.Range($bmPos, $bmPos + 20)
$cleanupRange.Text = ""

# ---------------------------------------------------------------------------
# 5) Styles: add the "Hyperlink" and "Mention" character styles referenced by
#    the new runs (best effort; some cosmetic attributes are not reachable
#    through the COM surface exposed by this host).
# ---------------------------------------------------------------------------
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = $d.Styles.Item("Default Paragraph Font")
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Color = 0xC16305
$hlStyle.Font.Underline = 1

$mentionStyle = $d.Styles.Add("Mention", 2)
$mentionStyle.BaseStyle = $d.Styles.Item("Default Paragraph Font")
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionStyle.Font.Color = 0x9A572B
